$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from G1 (existing "sum" header) onto the new H1 cell
# so the new "Save" header matches the look of the rest of the header row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header label and the data value for the new "Save" column
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
